# Updates the cryptocurrency listing on Sheet1 (rows 2-51) to reflect a
# newer data pull: refreshed "Price" (column D) and "Volume(1h)" (column E)
# figures for every coin, plus a reordering of two adjacent row pairs
# (FirstDigitalUSD/Algorand at rows 44-45, and Aave/Celestia at rows 46-47)
# whose Coin/Link/Price/Volume values swap places.
#
# Column D holds price text that looks numeric (e.g. "103.40", "0.0000105",
# "41.964.34"). Assigning such a string straight to Range.Value lets Excel's
# automatic type inference coerce it into a real number (dropping trailing
# zeros / using scientific notation), which corrupts the text. To keep these
# as genuine text values (matching the original inlineStr cells) we briefly
# force the cell to Text format, assign the string, then restore the
# cell's original ("Normal") style so no visible formatting changes stick.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "41.964.34"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -2.28%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.292.52"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "316.93"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "103.40"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -4.03%  "
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -3.11%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "39.30"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.86%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0905"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.72%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.26"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("E13").Value = "  -0.43%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.959"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.96%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.21"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.86%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.640.75"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.84%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.285.06"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.48%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "41.970.36"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.17%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.30%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0000105"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  -0.93%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "73.28"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.08%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "278.69"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.14%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.20"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +7.67%  "
$ws.Range("E25").Value = "  -2.79%  "
$ws.Range("E26").Value = "  +0.59%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.77"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -6.14%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.81%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "22.78"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.15%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "35.96"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.37%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "163.17"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.96%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0871"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("E36").Value = "  -5.03%  "
$ws.Range("E37").Value = "  -5.43%  "
$ws.Range("E38").Value = "  -4.57%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.83"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("E40").Value = "  -3.55%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "98.99"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -6.14%  "
$ws.Range("E42").Value = "  -5.05%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "69.30"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.224"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -6.65%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "11.90"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "112.46"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.41%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "77.03"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.89"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.04%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "5.28"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.99%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.581.34"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
